$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.071.98"
$ws.Range("E2").Value = "  -0.26%  "
$ws.Range("D3").Value = "2.337.07"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("E4").Value = "  +0.00%  "
$cell = $ws.Range("D5")
$cell.Value = "'303.49"
$cell.Style = "Normal"
$cell = $ws.Range("D6")
$cell.Value = "'94.68"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -3.25%  "
$cell = $ws.Range("D7")
$cell.Value = "'0.505"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -0.50%  "
$ws.Range("E8").Value = "  -0.01%  "
$cell = $ws.Range("D9")
$cell.Value = "'0.497"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -1.34%  "
$cell = $ws.Range("D10")
$cell.Value = "'34.21"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  -3.99%  "
$cell = $ws.Range("D11")
$cell.Value = "'19.00"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.66%  "
$ws.Range("E12").Value = "  -0.66%  "
$ws.Range("E13").Value = "  +2.39%  "
$cell = $ws.Range("D14")
$cell.Value = "'6.71"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -3.09%  "
$ws.Range("D15").Value = "2.700.43"
$ws.Range("E15").Value = "  +0.33%  "
$ws.Range("D16").Value = "2.364.52"
$ws.Range("E16").Value = "  +2.02%  "
$cell = $ws.Range("D17")
$cell.Value = "'0.792"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").Value = "43.010.76"
$ws.Range("E18").Value = "  -0.20%  "
$cell = $ws.Range("D19")
$cell.Value = "'12.11"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -3.53%  "
$cell = $ws.Range("D20")
$cell.Value = "'6.17"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").Value = "0.0₃0892"
$ws.Range("E21").Value = "  -0.74%  "
$cell = $ws.Range("D22")
$cell.Value = "'68.10"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +0.16%  "
$cell = $ws.Range("D23")
$cell.Value = "'236.82"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  -0.47%  "
$cell = $ws.Range("D24")
$cell.Value = "'2.24"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +1.08%  "
$ws.Range("E25").Value = "  -0.01%  "
$ws.Range("E26").Value = "  -0.63%  "
$cell = $ws.Range("D27")
$cell.Value = "'24.70"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -1.22%  "
$cell = $ws.Range("D28")
$cell.Value = "'2.05"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -14.13%  "
$cell = $ws.Range("D29")
$cell.Value = "'9.13"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +0.16%  "
$cell = $ws.Range("D30")
$cell.Value = "'31.68"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -4.41%  "
$ws.Range("B31").Value = "FirstDigitalUSD"
$ws.Range("C31").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$cell = $ws.Range("D31")
$cell.Value = "'1.00"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  +0.05%  "
$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$cell = $ws.Range("D32")
$cell.Value = "'140.06"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -15.74%  "
$ws.Range("E33").Value = "  +0.08%  "
$cell = $ws.Range("D34")
$cell.Value = "'17.72"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -1.85%  "
$ws.Range("E35").Value = "  +1.01%  "
$cell = $ws.Range("D36")
$cell.Value = "'4.39"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -2.70%  "
$cell = $ws.Range("D37")
$cell.Value = "'1.82"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("E38").Value = "  -1.55%  "
$ws.Range("E39").Value = "  -0.56%  "
$cell = $ws.Range("D40")
$cell.Value = "'22.37"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +24.12%  "
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("E42").Value = "  -0.65%  "
$ws.Range("D43").Value = "1.944.44"
$ws.Range("E43").Value = "  -2.82%  "
$cell = $ws.Range("D44")
$cell.Value = "'0.0280"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -0.60%  "
$cell = $ws.Range("D45")
$cell.Value = "'10.19"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -5.03%  "
$ws.Range("E46").Value = "  -2.20%  "
$cell = $ws.Range("D47")
$cell.Value = "'2.73"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.86%  "
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").Value = "2.565.00"
$ws.Range("E49").Value = "  +0.28%  "
$cell = $ws.Range("D50")
$cell.Value = "'52.91"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.58%  "
$cell = $ws.Range("D51")
$cell.Value = "'72.56"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  +0.62%  "
